# v1.7: Solapa Errores - columna Mes-Año reemplazada por Id_Origen
# Adds one new bitacora entry row to "Log" and one new version row to "Versiones".

$wb = $excel.ActiveWorkbook

# --- Sheet "Log": append row 38 describing the change ---
$log = $wb.Worksheets.Item("Log")

$log.Range("A38").Value = "27/02/2025"
$log.Range("B38").Value = "19:50"
$log.Range("C38").Value = "Solapa Errores: columna Mes-Año por Id_Origen"
$log.Range("D38").Value = "En la tabla de la solapa Errores se reemplaza la columna Mes-Año por Id_Origen (identificador de origen del registro)."
$log.Range("E38").Value = "Diagnostico"

# --- Sheet "Versiones": append row 9 for version 1.7 ---
$versiones = $wb.Worksheets.Item("Versiones")

# Column A holds version numbers stored as text (e.g. "1.0" .. "1.6"); force
# text formatting first so "1.7" isn't auto-coerced into a numeric value.
$versiones.Range("A9").NumberFormat = "@"
$versiones.Range("A9").Value = "1.7"
$versiones.Range("B9").Value = "27/02/2025"
$versiones.Range("C9").Value = "Solapa Errores: columna Mes-Año reemplazada por Id_Origen en la tabla"
